$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Fix the absolute/relative row references in the existing formulas
#    (the $6 anchors on C/E/G/U columns lose their row-dollar, W16 gains one)
# ---------------------------------------------------------------------------
$ws.Range("D8").Formula  = '=(C6-C8)*C7'
$ws.Range("F8").Formula  = '=(E6-E8)*E7'
$ws.Range("H8").Formula  = '=(G6-G8)*G7'
$ws.Range("V8").Formula  = '=(U6-U8)*U7'
$ws.Range("D10").Formula = '=(C6-C10)*C9'
$ws.Range("F10").Formula = '=(E6-E10)*E9'
$ws.Range("H10").Formula = '=(G6-G10)*G9'
$ws.Range("V10").Formula = '=(U6-U10)*U9'
$ws.Range("W16").Formula = '=W15*$E$1'

# ---------------------------------------------------------------------------
# 2. New "SAMPLE" header row (26) - four merged, centered header cells
# ---------------------------------------------------------------------------
$ws.Range("C26").HorizontalAlignment = -4108
$ws.Range("C26:D26").Merge()
$ws.Range("C26").Value = "BUY ORDER SAMPLE"

$ws.Range("E26").HorizontalAlignment = -4108
$ws.Range("E26:F26").Merge()
$ws.Range("E26").Value = "SELL ORDER SAMPLE"

$ws.Range("G26").HorizontalAlignment = -4108
$ws.Range("G26:H26").Merge()
$ws.Range("G26").Value = "SHORT ORDER SAMPLE"

$ws.Range("I26").HorizontalAlignment = -4108
$ws.Range("I26:J26").Merge()
$ws.Range("I26").Value = "CLOSE_SHORT SAMPLE"

# ---------------------------------------------------------------------------
# 3. Column headers (row 28)
# ---------------------------------------------------------------------------
$ws.Range("C28").Value = "BUY"
$ws.Range("D28").Value = "balance change"
$ws.Range("E28").Value = "SELL"
$ws.Range("F28").Value = "balance change"
$ws.Range("G28").Value = "SHORT"
$ws.Range("H28").Value = "balance change"
$ws.Range("I28").Value = "CLOSE_SHORT"

# ---------------------------------------------------------------------------
# 4. Sample order-book data (rows 29-36)
# ---------------------------------------------------------------------------
$ws.Range("A29").Value = "QTY"
$ws.Range("C29").Value = 120
$ws.Range("E29").Value = -100
$ws.Range("G29").Value = -100

$ws.Range("A30").Value = "PRICE"
$ws.Range("C30").Value = 10
$ws.Range("D30").Formula = '=-(C29*C30)'
$ws.Range("E30").Value = 10
$ws.Range("F30").Value = 0
$ws.Range("G30").Value = 11
$ws.Range("H30").Formula = '=$E$1*G29'

$ws.Range("A31").Value = "FILL QTY"
$ws.Range("C31").Value = 100
$ws.Range("E31").Value = 50
$ws.Range("G31").Value = 50

$ws.Range("A32").Value = "FILL PRICE"
$ws.Range("C32").Value = 9
$ws.Range("D32").Formula = '=(C30-C32)*C31'
$ws.Range("E32").Value = 11
$ws.Range("F32").Formula = '=E32*E31'
$ws.Range("G32").Value = 11
$ws.Range("H32").Formula = '=G32*G31'

$ws.Range("A33").Value = "FILL QTY"
$ws.Range("C33").Value = 20
$ws.Range("E33").Value = 10
$ws.Range("G33").Value = 10

$ws.Range("A34").Value = "FILL PRICE"
$ws.Range("C34").Value = 10
$ws.Range("D34").Formula = '=(C30-C34)*C33'
$ws.Range("E34").Value = 10
$ws.Range("F34").Formula = '=E34*E33'
$ws.Range("G34").Value = 12
$ws.Range("H34").Formula = '=G34*G33'

$ws.Range("A35").Value = "CANCEL QTY"
$ws.Range("C35").Value = 0
$ws.Range("D35").Formula = '=C35*C30'
$ws.Range("G35").Value = 0

$ws.Range("A36").Value = "OPEN QTY"
$ws.Range("C36").Formula = '=ABS(C29)-C31-C33-C35'
$ws.Range("E36").Formula = '=ABS(E29)-E31-E33-E35'
$ws.Range("G36").Formula = '=ABS(G29)-G31-G33-G35'

# ---------------------------------------------------------------------------
# 5. Summary rows (38-40)
# ---------------------------------------------------------------------------
$ws.Range("A38").Value = "POSITIVE POSITION"
$ws.Range("B38").Formula = '=SUM(C38:QT38)'
$ws.Range("C38").Formula = '=C31+C33'
$ws.Range("E38").Formula = '=-(E31+E33)'
$ws.Range("I38").Value = -60

$ws.Range("A39").Value = "SHORT POSITION"
$ws.Range("B39").Formula = '=SUM(C39:QT39)'
$ws.Range("G39").Formula = '=-(G31+G33)'
$ws.Range("I39").Value = 60

$ws.Range("A40").Value = "BALANCE CHANGE"
$ws.Range("B40").Formula = '=SUM(C40:CF40)'
$ws.Range("D40").Formula = '=SUM(D29:D39)'
$ws.Range("F40").Formula = '=SUM(F29:F39)'
$ws.Range("H40").Formula = '=SUM(H29:H39)'
$ws.Range("I40").Formula = '=I39*$E$1'

# ---------------------------------------------------------------------------
# 6. Move the active selection to match the final cursor position
# ---------------------------------------------------------------------------
$ws.Range("L25").Select()
